$d = $word.ActiveDocument

function Add-ParaAfterLast() {
    $count = $d.Paragraphs.Count
    $lastPara = $d.Paragraphs.Item($count)
    $r = $lastPara.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $newCount = $d.Paragraphs.Count
    return $d.Paragraphs.Item($newCount)
}

# --- 3 blank paragraphs (firstLine inherited from the preceding paragraph: 708 twips = 35.4pt) ---
Add-ParaAfterLast | Out-Null
Add-ParaAfterLast | Out-Null
Add-ParaAfterLast | Out-Null

# --- "3) SCRUM" heading paragraph (firstLine = 0) ---
$p = Add-ParaAfterLast
$p.Range.Text = "3) SCRUM"
$p.Range.ParagraphFormat.FirstLineIndent = 0

# --- body paragraph 1 (firstLine = 708 twips = 35.4pt) ---
$p = Add-ParaAfterLast
$p.Range.Text = "Métodos ágeis como o SCRUM são ótimas ferramentas a serem aplicadas a projetos em espiral, pois são totalmente compatíveis."
$p.Range.ParagraphFormat.FirstLineIndent = 35.4

# --- body paragraph 2 ---
$p = Add-ParaAfterLast
$p.Range.Text = "Os métodos ágeis consistem em administrar da melhor forma possível a equipe e o tempo. A forma de trabalho se dá pela divisão de tarefas e alinhamento constante do que cada membro está fazendo para entregar o produto. "

# --- body paragraph 3 ---
$p = Add-ParaAfterLast
$p.Range.Text = "Ou seja, combina totalmente com o modelo espiral, pois auxilia na modularização do projeto, na administração de tarefas, no controle de entregas e no alinhamento da equipe como um todo."

Write-Host "Final paragraph count: $($d.Paragraphs.Count)"
